$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 7238
$ws.Range("F4").Value = 3500
$ws.Range("F6").Value = 3830
$ws.Range("F7").Value = 64
$ws.Range("F8").Value = 72
$ws.Range("F10").Value = 96
$ws.Range("F11").Value = 138
$ws.Range("F14").Value = 126
$ws.Range("F15").Value = 361
$ws.Range("F18").Value = 350
$ws.Range("F19").Value = 4097
$ws.Range("F22").Value = 1025
$ws.Range("F23").Value = 531
$ws.Range("F24").Value = 1641
$ws.Range("F26").Value = 94
$ws.Range("F27").Value = 2999
$ws.Range("F28").Value = 2193
$ws.Range("F29").Value = 59
$ws.Range("F32").Value = 20
$ws.Range("F33").Value = 83
$ws.Range("F34").Value = 40
$ws.Range("F36").Value = 4264
$ws.Range("F37").Value = 460
$ws.Range("F38").Value = 320
$ws.Range("F41").Value = 781
$ws.Range("F42").Value = 193
$ws.Range("F43").Value = 10
$ws.Range("F44").Value = 1620
$ws.Range("F46").Value = 28
$ws.Range("F47").Value = 600
$ws.Range("F48").Value = 709

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 249
$ws.Range("F6").Value = 64
$ws.Range("F16").Value = 562

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 163

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 163
$ws.Range("F3").Value = 249
$ws.Range("F5").Value = 7238
$ws.Range("F6").Value = 3500
$ws.Range("F7").Value = 3500
$ws.Range("F8").Value = 3830
$ws.Range("F9").Value = 72
$ws.Range("F11").Value = 96
$ws.Range("F13").Value = 138
$ws.Range("F15").Value = 64
$ws.Range("F16").Value = 126
$ws.Range("F17").Value = 361
$ws.Range("F20").Value = 350
$ws.Range("F21").Value = 4097
$ws.Range("F26").Value = 1025
$ws.Range("F27").Value = 531
$ws.Range("F28").Value = 1641
$ws.Range("F30").Value = 94
$ws.Range("F31").Value = 3000
$ws.Range("F32").Value = 2193
$ws.Range("F33").Value = 59
$ws.Range("F36").Value = 40
$ws.Range("F39").Value = 4264
$ws.Range("F41").Value = 460
$ws.Range("F42").Value = 320
$ws.Range("F45").Value = 781
$ws.Range("F46").Value = 193
$ws.Range("F47").Value = 1620
$ws.Range("F49").Value = 600
$ws.Range("F50").Value = 709
